# "pruebas/pantalla inicio.xlsx" update
# - Sheet "Pruebas": update project/screen header labels and a few test-case
#   descriptions to proper sentence case / clearer wording, and move the
#   selection to F4 (reflecting the last edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pruebas")

$ws.Range("A1").Value = "Proyecto BlueWeb"
$ws.Range("A2").Value = "Pantalla Inicio"

$ws.Range("B4").Value = "Mostrar galeria"
$ws.Range("E4").Value = "Cargar la pantalla de inicio donde se muestra la galeria"
$ws.Range("F4").Value = "Muestra la galeria de imágenes"

$ws.Range("B5").Value = "Funcionalidad de botones "
$ws.Range("E5").Value = "Hacer click en cada boton del menu "
$ws.Range("F5").Value = "Redirije a las paginas con éxito "

[void]$ws.Range("F4").Select()
